$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data: Row -> (D, J, K, L, M, P)
$data = @{
    2 = @(44396, 130, 22000, 22000, 22000, 1467)
    3 = @(44446, 150, 22000, 24000, 22667, 1511)
    4 = @(44400, 130, 24000, 24000, 24000, 1600)
    5 = @(44398, 130, 20000, 20000, 20000, 1333)
    6 = @(44391, 160, 20000, 20000, 20000, 1333)
    7 = @(44435, 140, 21000, 23000, 21714, 1448)
    8 = @(44483, 220, 18000, 20000, 18909, 1261)
    9 = @(44406, 400, 20000, 22000, 20850, 1390)
    10 = @(44365, 580, 20000, 22000, 21103, 1407)
    11 = @(44399, 150, 22000, 22000, 22000, 1467)
    12 = @(44476, 220, 20000, 22000, 20909, 1394)
    13 = @(44392, 220, 23000, 23000, 23000, 1533)
    14 = @(44453, 280, 20000, 22000, 21286, 1419)
    15 = @(44449, 220, 22000, 24000, 23091, 1539)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]  # P: Precio $/Kg
}
